$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.879.39'
$ws.Range("E2").Value = '  +0.59%  '

$ws.Range("D3").Value = '3.358.66'
$ws.Range("E3").Value = '  +0.27%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '554.80'
$ws.Range("E5").Value = '  -0.13%  '

$ws.Range("D6").Value = '174.04'
$ws.Range("E6").Value = '  -0.96%  '

$ws.Range("E7").Value = '  +1.99%  '

$ws.Range("D8").Value = '3.350.00'
$ws.Range("E8").Value = '  +0.27%  '

$ws.Range("E10").Value = '  +6.59%  '

$ws.Range("E11").Value = '  +1.40%  '

$ws.Range("D12").Value = '53.69'
$ws.Range("E12").Value = '  -1.99%  '

$ws.Range("E13").Value = '  +3.11%  '

$ws.Range("D14").Value = '9.13'
$ws.Range("E14").Value = '  +0.91%  '

$ws.Range("D15").Value = '3.887.57'
$ws.Range("E15").Value = '  +0.07%  '

$ws.Range("E16").Value = '  +2.31%  '

$ws.Range("D17").Value = '18.21'
$ws.Range("E17").Value = '  -0.53%  '

$ws.Range("D18").Value = '3.345.12'
$ws.Range("E18").Value = '  -0.02%  '

$ws.Range("D19").Value = '64.761.23'
$ws.Range("E19").Value = '  +0.57%  '

$ws.Range("D20").Value = '11.80'
$ws.Range("E20").Value = '  +0.37%  '

$ws.Range("D21").Value = '0.992'
$ws.Range("E21").Value = '  +1.38%  '

$ws.Range("D22").Value = '449.21'
$ws.Range("E22").Value = '  +3.06%  '

$ws.Range("E23").Value = '  -1.78%  '

$ws.Range("D24").Value = '4.06'
$ws.Range("E24").Value = '  -0.32%  '

$ws.Range("D25").Value = '86.75'
$ws.Range("E25").Value = '  +2.85%  '

$ws.Range("D26").Value = '13.71'
$ws.Range("E26").Value = '  +2.29%  '

$ws.Range("D27").Value = '2.87'
$ws.Range("E27").Value = '  +1.19%  '

$ws.Range("D28").Value = '10.73'
$ws.Range("E28").Value = '  -0.20%  '

$ws.Range("D29").Value = '8.64'
$ws.Range("E29").Value = '  -0.98%  '

$ws.Range("D30").Value = '30.96'
$ws.Range("E30").Value = '  +4.10%  '

$ws.Range("D31").Value = '6.55'
$ws.Range("E31").Value = '  -1.31%  '

$ws.Range("D32").Value = '63.06'
$ws.Range("E32").Value = '  +7.84%  '

$ws.Range("D33").Value = '11.43'
$ws.Range("E33").Value = '  -0.38%  '

$ws.Range("D34").Value = '576.52'
$ws.Range("E34").Value = '  -0.24%  '

$ws.Range("E35").Value = '  -0.20%  '

$ws.Range("E36").Value = '  +0.01%  '

$ws.Range("D37").Value = '3.64'
$ws.Range("E37").Value = '  +4.06%  '

$ws.Range("E38").Value = '  -0.35%  '

$ws.Range("D39").Value = '35.66'
$ws.Range("E39").Value = '  -0.20%  '

$ws.Range("E40").Value = '  +0.78%  '

$ws.Range("D41").Value = '0.0₃0741'
$ws.Range("E41").Value = '  -1.99%  '

$ws.Range("D42").Value = '3.077.16'
$ws.Range("E42").Value = '  -0.84%  '

$ws.Range("E43").Value = '  +1.53%  '

$ws.Range("D44").Value = '2.75'
$ws.Range("E44").Value = '  -2.08%  '

$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").Value = '3.20'
$ws.Range("E45").Value = '  -1.36%  '

$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").Value = '0.134'
$ws.Range("E46").Value = '  +3.37%  '

$ws.Range("E47").Value = '  -1.09%  '

$ws.Range("D48").Value = '0.998'

$ws.Range("D49").Value = '141.47'
$ws.Range("E49").Value = '  +4.23%  '

$ws.Range("E50").Value = '  -2.71%  '

$ws.Range("D51").Value = '8.27'
$ws.Range("E51").Value = '  +0.06%  '
